# 自动更新Excel文件: daily roll-over of the "剩余" (days remaining) counter.
#
# For every data row (row 2 .. last used row):
#   D = 总天 (total days for this delivery cycle)
#   E = 剩余 (days remaining)
#   F = 开始时间 (cycle start date, yyyymmdd)
#
# Each day:
#   - if the cycle was just (re)started today (E already equals D), leave it alone
#   - if the cycle has run out (E <= 1), refill it: reset E back to the full
#     cycle length (D) and stamp F with today's date
#   - otherwise simply count the day down: E = E - 1

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$today = 20251211

$lastRow = $ws.Cells.Item(1, 1).End(4).Row

for ($r = 2; $r -le $lastRow; $r++) {
    $dCell = $ws.Cells.Item($r, 4)
    $eCell = $ws.Cells.Item($r, 5)
    $fCell = $ws.Cells.Item($r, 6)

    $d = $dCell.Value()
    $e = $eCell.Value()

    if ($e -eq $null -or $d -eq $null) {
        continue
    }

    if ($e -eq $d) {
        # cycle just started (e.g. refilled earlier today) - nothing to do
        continue
    }
    elseif ($e -le 1) {
        # out of days - restock/refill, cycle restarts today
        $eCell.Value = $d
        $fCell.Value = $today
    }
    else {
        # one more day has passed
        $eCell.Value = $e - 1
    }
}
